$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2142.125
$ws.Range("I2").Value = 2687.4
$ws.Range("J2").Value = 1233.3334
$ws.Range("K2").Value = 2687.4
$ws.Range("L2").Value = 1233.3334
$ws.Range("M2").Value = -2574.4
$ws.Range("N2").Value = -1459.3334
$ws.Range("H62").Value = 7992.143
$ws.Range("J62").Value = 7656.6665
$ws.Range("L62").Value = 7656.6665
$ws.Range("N62").Value = -8904.666499999999
$ws.Range("H65").Value = 7992.143
$ws.Range("J65").Value = 7656.6665
$ws.Range("L65").Value = 38283.3325
$ws.Range("N65").Value = -44523.3325
$ws.Range("H125").Value = 2669.889
$ws.Range("I125").Value = 1532
$ws.Range("J125").Value = 2812.125
$ws.Range("K125").Value = 13788
$ws.Range("L125").Value = 25309.125
$ws.Range("M125").Value = -11328
$ws.Range("N125").Value = -30229.125
$ws.Range("H132").Value = 2090.3872
$ws.Range("J132").Value = 6725
$ws.Range("L132").Value = 20175
$ws.Range("N132").Value = -25235

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1300
$ws.Range("I22").Value = 1300
$ws.Range("K22").Value = 1300
$ws.Range("M22").Value = -1001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1963.7222
$ws.Range("I105").Value = 1379.0952
$ws.Range("J105").Value = 2782.2
$ws.Range("K105").Value = 1379.0952
$ws.Range("L105").Value = 2782.2
$ws.Range("M105").Value = 367.9048
$ws.Range("N105").Value = -6276.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 4596
$ws.Range("I7").Value = 246.38889
$ws.Range("J7").Value = 20254.6
$ws.Range("K7").Value = 246.38889
$ws.Range("L7").Value = 20254.6
$ws.Range("M7").Value = -133.38889
$ws.Range("N7").Value = -20480.6
$ws.Range("H26").Value = 14370.37
$ws.Range("J26").Value = 14370.37
$ws.Range("L26").Value = 14370.37
$ws.Range("N26").Value = -14944.37
$ws.Range("H31").Value = 1029223.94
$ws.Range("J31").Value = 1306547.6
$ws.Range("L31").Value = 1306547.6
$ws.Range("N31").Value = -1307137.6
$ws.Range("H34").Value = 1029223.94
$ws.Range("J34").Value = 1306547.6
$ws.Range("L34").Value = 1306547.6
$ws.Range("N34").Value = -1306951.6
$ws.Range("H58").Value = 4294.5
$ws.Range("I58").Value = 4532.6665
$ws.Range("K58").Value = 4532.6665
$ws.Range("M58").Value = -4329.6665
$ws.Range("H93").Value = 22286.715
$ws.Range("I93").Value = 8334.5
$ws.Range("J93").Value = 106000
$ws.Range("K93").Value = 8334.5
$ws.Range("L93").Value = 106000
$ws.Range("M93").Value = -6462.5
$ws.Range("N93").Value = -109744
$ws.Range("H105").Value = 1790.5385
$ws.Range("I105").Value = 2029.6666
$ws.Range("K105").Value = 2029.6666
$ws.Range("M105").Value = -282.6666
$ws.Range("H134").Value = 3102.1482
$ws.Range("I134").Value = 1957.5333
$ws.Range("J134").Value = 4532.9165
$ws.Range("K134").Value = 5872.5999
$ws.Range("L134").Value = 13598.7495
$ws.Range("M134").Value = -3337.5999
$ws.Range("N134").Value = -18668.7495
$ws.Range("H136").Value = 4294.5
$ws.Range("I136").Value = 4532.6665
$ws.Range("K136").Value = 13597.9995
$ws.Range("M136").Value = -11047.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 205.04878
$ws.Range("I2").Value = 51.1
$ws.Range("K2").Value = 306.6
$ws.Range("M2").Value = -193.6
$ws.Range("H39").Value = 3472.9092
$ws.Range("J39").Value = 2311.75
$ws.Range("L39").Value = 6935.25
$ws.Range("N39").Value = -7523.25
$ws.Range("H50").Value = 412.8421
$ws.Range("I50").Value = 274.16666
$ws.Range("J50").Value = 476.84616
$ws.Range("K50").Value = 822.4999799999999
$ws.Range("L50").Value = 1430.53848
$ws.Range("M50").Value = -341.4999799999999
$ws.Range("N50").Value = -2392.53848
$ws.Range("H53").Value = 412.8421
$ws.Range("I53").Value = 274.16666
$ws.Range("J53").Value = 476.84616
$ws.Range("K53").Value = 822.4999799999999
$ws.Range("L53").Value = 1430.53848
$ws.Range("M53").Value = -341.4999799999999
$ws.Range("N53").Value = -2392.53848
$ws.Range("H68").Value = 2863.3
$ws.Range("I68").Value = 1916.6666
$ws.Range("K68").Value = 5749.9998
$ws.Range("M68").Value = -4938.9998
$ws.Range("H71").Value = 2863.3
$ws.Range("I71").Value = 1916.6666
$ws.Range("K71").Value = 17249.9994
$ws.Range("M71").Value = -13193.9994
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H132").Value = 2518.4375
$ws.Range("J132").Value = 2666.3333
$ws.Range("L132").Value = 23996.9997
$ws.Range("N132").Value = -29056.9997
$ws.Range("H133").Value = 2628.2856
$ws.Range("I133").Value = 2233.1667
$ws.Range("J133").Value = 4999
$ws.Range("K133").Value = 6699.500100000001
$ws.Range("L133").Value = 14997
$ws.Range("M133").Value = -1639.500100000001
$ws.Range("N133").Value = -25117
$ws.Range("H134").Value = 6065.3213
$ws.Range("I134").Value = 1578.8235
$ws.Range("K134").Value = 4736.470499999999
$ws.Range("M134").Value = 333.5295000000006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 32999.6
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H102").Value = 2520.353
$ws.Range("I102").Value = 1528.2941
$ws.Range("J102").Value = 3512.4119
$ws.Range("K102").Value = 1528.2941
$ws.Range("L102").Value = 3512.4119
$ws.Range("M102").Value = 93.70589999999993
$ws.Range("N102").Value = -6756.4119
$ws.Range("H126").Value = 2758.92
$ws.Range("I126").Value = 1887.3889
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 5662.1667
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -3192.1667
$ws.Range("N126").Value = -19940
$ws.Range("H136").Value = 8763.85
$ws.Range("J136").Value = 8763.85
$ws.Range("L136").Value = 26291.55
$ws.Range("N136").Value = -31391.55
$ws.Range("H139").Value = 49799.8
$ws.Range("J139").Value = 49799.8
$ws.Range("L139").Value = 49799.8
$ws.Range("N139").Value = -60079.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 486624.9
$ws.Range("I7").Value = 7460.9375
$ws.Range("K7").Value = 7460.9375
$ws.Range("M7").Value = -7348.9375
$ws.Range("H22").Value = 1366.3
$ws.Range("I22").Value = 884.8
$ws.Range("J22").Value = 1847.8
$ws.Range("K22").Value = 884.8
$ws.Range("L22").Value = 1847.8
$ws.Range("M22").Value = -589.8
$ws.Range("N22").Value = -2437.8
$ws.Range("H27").Value = 1366.3
$ws.Range("I27").Value = 884.8
$ws.Range("J27").Value = 1847.8
$ws.Range("K27").Value = 884.8
$ws.Range("L27").Value = 1847.8
$ws.Range("M27").Value = -777.8
$ws.Range("N27").Value = -2061.8
$ws.Range("H40").Value = 3851.0715
$ws.Range("I40").Value = 3346.6
$ws.Range("J40").Value = 5112.25
$ws.Range("K40").Value = 3346.6
$ws.Range("L40").Value = 5112.25
$ws.Range("M40").Value = -3210.6
$ws.Range("N40").Value = -5384.25
$ws.Range("H61").Value = 1802.5454
$ws.Range("I61").Value = 1261.2858
$ws.Range("J61").Value = 2749.75
$ws.Range("K61").Value = 1261.2858
$ws.Range("L61").Value = 2749.75
$ws.Range("M61").Value = -1059.2858
$ws.Range("N61").Value = -3153.75
$ws.Range("H63").Value = 111000
$ws.Range("J63").Value = 111000
$ws.Range("L63").Value = 111000
$ws.Range("N63").Value = -112498
$ws.Range("H66").Value = 111000
$ws.Range("J66").Value = 111000
$ws.Range("L66").Value = 333000
$ws.Range("N66").Value = -340488
$ws.Range("H113").Value = 1802.5454
$ws.Range("I113").Value = 1261.2858
$ws.Range("J113").Value = 2749.75
$ws.Range("K113").Value = 1261.2858
$ws.Range("L113").Value = 2749.75
$ws.Range("M113").Value = 908.7141999999999
$ws.Range("N113").Value = -7089.75
$ws.Range("H126").Value = 486624.9
$ws.Range("I126").Value = 7460.9375
$ws.Range("K126").Value = 22382.8125
$ws.Range("M126").Value = -19912.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 32247.5
$ws.Range("I47").Value = 27247.5
$ws.Range("K47").Value = 27247.5
$ws.Range("M47").Value = -26675.5
$ws.Range("H70").Value = 112999
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 112999
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 112999
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -113629
$ws.Range("H73").Value = 112999
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 112999
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 112999
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -115183
$ws.Range("H87").Value = 55000
$ws.Range("J87").Value = 55000
$ws.Range("L87").Value = 55000
$ws.Range("N87").Value = -57496
$ws.Range("H90").Value = 55000
$ws.Range("J90").Value = 55000
$ws.Range("L90").Value = 165000
$ws.Range("N90").Value = -177480
$ws.Range("H96").Value = 3336.2856
$ws.Range("I96").Value = 2417.1667
$ws.Range("J96").Value = 4025.625
$ws.Range("K96").Value = 2417.1667
$ws.Range("L96").Value = 4025.625
$ws.Range("M96").Value = -1044.1667
$ws.Range("N96").Value = -6771.625
$ws.Range("H126").Value = 1997.8125
$ws.Range("I126").Value = 2074.2307
$ws.Range("K126").Value = 6222.6921
$ws.Range("M126").Value = -3752.6921
